$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 33

# Row 3
$ws.Range("B3").Value = "<kilo>"
$ws.Range("C3").Value = 32

# Row 4
$ws.Range("B4").Value = "<now>"
$ws.Range("C4").Value = 36

# Row 5
$ws.Range("B5").Value = "<quebec>"
$ws.Range("C5").Value = 28

# Row 6
$ws.Range("B6").Value = "<not>"
$ws.Range("C6").Value = 32

# Row 7
$ws.Range("B7").Value = "<otin>"
$ws.Range("C7").Value = 32

# Row 8
$ws.Range("C8").Value = 32

# Row 9
$ws.Range("C9").Value = 36

# Row 10
$ws.Range("C10").Value = 28

# Row 11
$ws.Range("B11").Value = "<an>"
$ws.Range("C11").Value = 25

# Row 12
$ws.Range("C12").Value = 27

# Row 13
$ws.Range("C13").Value = 34

# Row 14
$ws.Range("B14").Value = "<sit>"
$ws.Range("C14").Value = 38

# Row 15
$ws.Range("B15").Value = "<it>"
$ws.Range("C15").Value = 27

# Row 16
$ws.Range("C16").Value = 36

# Row 17
$ws.Range("B17").Value = "<serbo>"
$ws.Range("C17").Value = 30

# Row 18
$ws.Range("C18").Value = 32
